$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "dSF" column (F) values per re-pulled data / mean calculation fix
$ws.Range("F5").Value = 6
$ws.Range("F8").Value = -8
$ws.Range("F9").Value = -3
$ws.Range("F19").Value = -7
$ws.Range("F20").Value = 1
